# Append: 2025-09-05 12:44 JST
# Update the "取得日時" (acquired datetime) timestamp in column A for all
# existing data rows (rows 2-15) from "2025-09-05 12:34:17" to
# "2025-09-05 12:44:18" on the active sheet ("ランサーズ").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2025-09-05 12:34:17"
$newTimestamp = "2025-09-05 12:44:18"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
